$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.65527733097282
$ws.Range("C2").Value = 9.217579372478697
$ws.Range("D2").Value = 5.992695946119151
$ws.Range("E2").Value = 11.45696887856955
$ws.Range("G2").Value = 3.648436921934586
$ws.Range("I2").Value = 24.04506748762702
$ws.Range("L2").Value = 10.09427547790971
$ws.Range("M2").Value = 15.13757940450871
$ws.Range("N2").Value = 18.3783084496608
$ws.Range("O2").Value = 24.87743058639047
$ws.Range("B3").Value = 15.2164130235177
$ws.Range("C3").Value = 8.878859553085061
$ws.Range("D3").Value = 5.875524423247644
$ws.Range("E3").Value = 11.48335211977613
$ws.Range("G3").Value = 3.650834671908555
$ws.Range("I3").Value = 24.13108588010045
$ws.Range("L3").Value = 10.10280143323905
$ws.Range("M3").Value = 15.05260877172884
$ws.Range("N3").Value = 18.4395106789159
$ws.Range("O3").Value = 24.90377155126863
$ws.Range("B4").Value = 14.94311038225234
$ws.Range("C4").Value = 8.662349398054443
$ws.Range("D4").Value = 5.80420264812263
$ws.Range("E4").Value = 11.50044092543329
$ws.Range("G4").Value = 3.652385608228251
$ws.Range("I4").Value = 24.18917336710144
$ws.Range("L4").Value = 10.1094802688458
$ws.Range("M4").Value = 15.00267361744759
$ws.Range("N4").Value = 18.4789237657561
$ws.Range("O4").Value = 24.92639514782073
$ws.Range("B5").Value = 14.83094552353581
$ws.Range("C5").Value = 8.572050925915702
$ws.Range("D5").Value = 5.775339402445063
$ws.Range("E5").Value = 11.5076289849975
$ws.Range("G5").Value = 3.653037483766122
$ws.Range("I5").Value = 24.21416679635724
$ws.Range("L5").Value = 10.1125654116447
$ws.Range("M5").Value = 14.98290314535522
$ws.Range("N5").Value = 18.49544758576067
$ws.Range("O5").Value = 24.93723283319563
$ws.Range("B6").Value = 14.81227807472335
$ws.Range("C6").Value = 8.556934426554973
$ws.Range("D6").Value = 5.77056012715784
$ws.Range("E6").Value = 11.50883611723761
$ws.Range("G6").Value = 3.653146928380158
$ws.Range("I6").Value = 24.21839670141433
$ws.Range("L6").Value = 10.11309965955112
$ws.Range("M6").Value = 14.97965566461621
$ws.Range("N6").Value = 18.49821933651466
$ws.Range("O6").Value = 24.939130050692
$ws.Range("B7").Value = 14.94160065727327
$ws.Range("C7").Value = 8.661139868783165
$ws.Range("D7").Value = 5.803812515539833
$ws.Range("E7").Value = 11.50053695742037
$ws.Range("G7").Value = 3.652394319165994
$ws.Range("I7").Value = 24.18950508775583
$ws.Range("L7").Value = 10.10952040404855
$ws.Range("M7").Value = 15.00240462311358
$ws.Range("N7").Value = 18.47914473670011
$ws.Range("O7").Value = 24.92653476114184
$ws.Range("B8").Value = 15.50485430725485
$ws.Range("C8").Value = 9.10260369712203
$ws.Range("D8").Value = 5.952192708571378
$ws.Range("E8").Value = 11.46588158236537
$ws.Range("G8").Value = 3.649247369188807
$ws.Range("I8").Value = 24.07363046798036
$ws.Range("L8").Value = 10.09691580778367
$ws.Range("M8").Value = 15.10782701646215
$ws.Range("N8").Value = 18.39903099090515
$ws.Range("O8").Value = 24.88517230605179
$ws.Range("B9").Value = 16.57186512369751
$ws.Range("C9").Value = 9.897656803618105
$ws.Range("D9").Value = 6.246171777909001
$ws.Range("E9").Value = 11.40495248855538
$ws.Range("G9").Value = 3.64369777921242
$ws.Range("I9").Value = 23.88837763957096
$ws.Range("L9").Value = 10.08363550454316
$ws.Range("M9").Value = 15.33157544516561
$ws.Range("N9").Value = 18.25642454359037
$ws.Range("O9").Value = 24.85536755914772
$ws.Range("B10").Value = 17.32424845644298
$ws.Range("C10").Value = 10.43536422090796
$ws.Range("D10").Value = 6.461476841036696
$ws.Range("E10").Value = 11.36443580609818
$ws.Range("G10").Value = 3.639995303431135
$ws.Range("I10").Value = 23.77806064019801
$ws.Range("L10").Value = 10.08082228052668
$ws.Range("M10").Value = 15.50530639693698
$ws.Range("N10").Value = 18.16040298159903
$ws.Range("O10").Value = 24.86487698995134
$ws.Range("B11").Value = 17.65807699556569
$ws.Range("C11").Value = 10.66934398301518
$ws.Range("D11").Value = 6.55877319091034
$ws.Range("E11").Value = 11.34691810232644
$ws.Range("G11").Value = 3.638391461396227
$ws.Range("I11").Value = 23.73351523871507
$ws.Range("L11").Value = 10.08104260006113
$ws.Range("M11").Value = 15.58613831409439
$ws.Range("N11").Value = 18.11860237407027
$ws.Range("O11").Value = 24.87603333795955
$ws.Range("B12").Value = 17.78315879308646
$ws.Range("C12").Value = 10.75637980525225
$ws.Range("D12").Value = 6.595484313716284
$ws.Range("E12").Value = 11.34041536968798
$ws.Range("G12").Value = 3.637795627691717
$ws.Range("I12").Value = 23.71746111023117
$ws.Range("L12").Value = 10.08134092732142
$ws.Range("M12").Value = 15.61698597403853
$ws.Range("N12").Value = 18.10304258886396
$ws.Range("O12").Value = 24.88123949935631
$ws.Range("B13").Value = 17.75628119727747
$ws.Range("C13").Value = 10.73770535167025
$ws.Range("D13").Value = 6.587584438837832
$ws.Range("E13").Value = 11.34181003858634
$ws.Range("G13").Value = 3.637923440332974
$ws.Range("I13").Value = 23.72088238612663
$ws.Range("L13").Value = 10.08126713344483
$ws.Range("M13").Value = 15.61033210864047
$ws.Range("N13").Value = 18.10638171258444
$ws.Range("O13").Value = 24.88007462808521
$ws.Range("B14").Value = 17.66839488326552
$ws.Range("C14").Value = 10.67653608684769
$ws.Range("D14").Value = 6.561796321138059
$ws.Range("E14").Value = 11.34638049932063
$ws.Range("G14").Value = 3.638342211532163
$ws.Range("I14").Value = 23.73217811770168
$ws.Range("L14").Value = 10.08106284172229
$ws.Range("M14").Value = 15.58867151159723
$ws.Range("N14").Value = 18.11731687261349
$ws.Range("O14").Value = 24.87644198868079
$ws.Range("B15").Value = 17.61438519812144
$ws.Range("C15").Value = 10.6388629781553
$ws.Range("D15").Value = 6.545981901368596
$ws.Range("E15").Value = 11.34919706246539
$ws.Range("G15").Value = 3.638600217701803
$ws.Range("I15").Value = 23.73920321955206
$ws.Range("L15").Value = 10.08096566752674
$ws.Range("M15").Value = 15.5754341957609
$ws.Range("N15").Value = 18.12404999735517
$ws.Range("O15").Value = 24.87434467510336
$ws.Range("B16").Value = 17.30225182460983
$ws.Range("C16").Value = 10.41985578602721
$ws.Range("D16").Value = 6.455101821250146
$ws.Range("E16").Value = 11.36559896704043
$ws.Range("G16").Value = 3.640101731682022
$ws.Range("I16").Value = 23.78108554310544
$ws.Range("L16").Value = 10.08083800000215
$ws.Range("M16").Value = 15.50005832670333
$ws.Range("N16").Value = 18.16317248245898
$ws.Range("O16").Value = 24.86428530184356
$ws.Range("B17").Value = 17.10851957900674
$ws.Range("C17").Value = 10.28275063213985
$ws.Range("D17").Value = 6.399155100953252
$ws.Range("E17").Value = 11.37589460870229
$ws.Range("G17").Value = 3.641043420039851
$ws.Range("I17").Value = 23.80822560333674
$ws.Range("L17").Value = 10.08114333426293
$ws.Range("M17").Value = 15.45426437939335
$ws.Range("N17").Value = 18.18765356609483
$ws.Range("O17").Value = 24.85986342087168
$ws.Range("B18").Value = 16.99630240326304
$ws.Range("C18").Value = 10.20289315795436
$ws.Range("D18").Value = 6.366917383624158
$ws.Range("E18").Value = 11.38190240703777
$ws.Range("G18").Value = 3.641592628526931
$ws.Range("I18").Value = 23.82436635592311
$ws.Range("L18").Value = 10.08146020207613
$ws.Range("M18").Value = 15.4280956736436
$ws.Range("N18").Value = 18.20191146810504
$ws.Range("O18").Value = 24.85796308659003
$ws.Range("B19").Value = 16.95817615174484
$ws.Range("C19").Value = 10.17568463451487
$ws.Range("D19").Value = 6.355993423199649
$ws.Range("E19").Value = 11.38395133397783
$ws.Range("G19").Value = 3.64177988379487
$ws.Range("I19").Value = 23.82992234561969
$ws.Range("L19").Value = 10.08159177078307
$ws.Range("M19").Value = 15.41926534927104
$ws.Range("N19").Value = 18.20676938938061
$ws.Range("O19").Value = 24.85743011167859
$ws.Range("B20").Value = 17.12922510966634
$ws.Range("C20").Value = 10.29744932367577
$ws.Range("D20").Value = 6.40511708966832
$ws.Range("E20").Value = 11.37478972139485
$ws.Range("G20").Value = 3.640942392233243
$ws.Range("I20").Value = 23.80528156661532
$ws.Range("L20").Value = 10.08109621791111
$ws.Range("M20").Value = 15.45912169807047
$ws.Range("N20").Value = 18.18502919940068
$ws.Range("O20").Value = 24.86026759035394
$ws.Range("B21").Value = 17.69424626140075
$ws.Range("C21").Value = 10.69454580164147
$ws.Range("D21").Value = 6.56937483012578
$ws.Range("E21").Value = 11.34503449736081
$ws.Range("G21").Value = 3.638218896470807
$ws.Range("I21").Value = 23.72883816199269
$ws.Range("L21").Value = 10.08111702180421
$ws.Range("M21").Value = 15.59502745062658
$ws.Range("N21").Value = 18.1140976542168
$ws.Range("O21").Value = 24.87748235561647
$ws.Range("B22").Value = 18.05571005951032
$ws.Range("C22").Value = 10.94492191992959
$ws.Range("D22").Value = 6.675935204549758
$ws.Range("E22").Value = 11.32635014390375
$ws.Range("G22").Value = 3.636505975670207
$ws.Range("I22").Value = 23.68362557180734
$ws.Range("L22").Value = 10.08238276504181
$ws.Range("M22").Value = 15.68522932293818
$ws.Range("N22").Value = 18.06930835525846
$ws.Range("O22").Value = 24.89445313317376
$ws.Range("B23").Value = 17.86354105872029
$ws.Range("C23").Value = 10.8121401701661
$ws.Range("D23").Value = 6.619146962775643
$ws.Range("E23").Value = 11.33625274659454
$ws.Range("G23").Value = 3.637414079070794
$ws.Range("I23").Value = 23.70732086067372
$ws.Range("L23").Value = 10.08159293297156
$ws.Range("M23").Value = 15.63696754083901
$ws.Range("N23").Value = 18.09307010047973
$ws.Range("O23").Value = 24.88487260484978
$ws.Range("B24").Value = 17.11986674002597
$ws.Range("C24").Value = 10.29080726263232
$ws.Range("D24").Value = 6.402421899548702
$ws.Range("E24").Value = 11.37528896463296
$ws.Range("G24").Value = 3.640988042543349
$ws.Range("I24").Value = 23.80661089133126
$ws.Range("L24").Value = 10.08111707893366
$ws.Range("M24").Value = 15.45692520793082
$ws.Range("N24").Value = 18.18621510418536
$ws.Range("O24").Value = 24.86008286587191
$ws.Range("B25").Value = 16.28816879135542
$ws.Range("C25").Value = 9.690513276549561
$ws.Range("D25").Value = 6.166585568686013
$ws.Range("E25").Value = 11.42068673545672
$ws.Range("G25").Value = 3.645132974339123
$ws.Range("I25").Value = 23.93397877157326
$ws.Range("L25").Value = 10.08600637589678
$ws.Range("M25").Value = 15.33157544516561
$ws.Range("N25").Value = 18.29346038307748
$ws.Range("O25").Value = 24.85792118652575
